$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

$ws.Range("B2").Value = 1.029200236871507
$ws.Range("C2").Value = 0.8534801061452516
$ws.Range("D2").Value = 0.8779275743016757
$ws.Range("E2").Value = 1.529865779888268
$ws.Range("F2").Value = 0.8170539553072627

$ws.Range("B3").Value = 2.389019732596685
$ws.Range("C3").Value = 2.287274680662983
$ws.Range("D3").Value = 1.750458490607735
$ws.Range("E3").Value = 2.457450037569061
$ws.Range("F3").Value = 1.689195688397789

$ws.Range("B4").Value = 0.9153642748603353
$ws.Range("C4").Value = 1.140785701675978
$ws.Range("D4").Value = 1.173347807821229
$ws.Range("E4").Value = 1.763832234636872
$ws.Range("F4").Value = 1.068988913966481

$ws.Range("B5").Value = 3.812634251111108
$ws.Range("C5").Value = 1.996457612222223
$ws.Range("D5").Value = 1.841608207777778
$ws.Range("E5").Value = 2.716257127777777
$ws.Range("F5").Value = 1.594091214444444

$wb.Save()
